$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$grp = $s.Shapes.Item(2)
$tb = $grp.GroupItems.Item(2)
$tr = $tb.TextFrame.TextRange

$bottomPara = $tr.Paragraphs(7)
$topPara = $tr.Paragraphs(13)

$bottomPara.Runs(1).Text = "Top 10% locations: "
$topPara.Runs(1).Text = "Bottom 10% locations: "
